$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.474.45"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "3.200.07"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.17"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.51"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.199.19"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.64"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("D15").Value = "3.726.15"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "3.195.51"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "62.577.26"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.23"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  -4.05%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.95"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("E31").Value = "  -3.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.45"
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.45"
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.67"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "0.0₃0695"
$ws.Range("E38").Value = "  -8.46%  "
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "416.06"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.997.92"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("E44").Value = "  -5.21%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.40"
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.00"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.10"
$ws.Range("E51").Value = "  -0.09%  "
